# Re-sequence observation rows 47-65 on the "Artfynd" sheet: the underlying
# source export renumbered each record's id (col A) and shuffled the full
# row content (species, activity, location, comments, ...) between rows.
# Below, each row is rewritten cell-by-cell (A:AY) to its final content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 47 ----
$ws.Range("A47").Value = 111736268
$ws.Range("B47").Value = 89401
$ws.Range("C47").Value = 'Ovaliderad'
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 1108
$ws.Range("F47").Value = 'Harticka'
$ws.Range("G47").Value = 'Pelloporus leporinus'
$ws.Range("H47").Value = '(Fr.) Krieglst.'
$ws.Range("I47").Value = ""
$ws.Range("J47").Value = ""
$ws.Range("K47").Value = ""
$ws.Range("L47").Value = ""
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = ""
$ws.Range("O47").Value = ""
$ws.Range("P47").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q47").Value = 616308.8236423519
$ws.Range("R47").Value = 7268903.133137755
$ws.Range("S47").Value = 10
$ws.Range("T47").Value = 'Västerbotten'
$ws.Range("U47").Value = 'Sorsele'
$ws.Range("V47").Value = 'Lycksele lappmark'
$ws.Range("W47").Value = 'Sorsele'
$ws.Range("X47").Value = ""
$ws.Range("Y47").NumberFormat = "@"
$ws.Range("Y47").Value = '2023-08-27'
$ws.Range("Z47").Value = '00:00'
$ws.Range("AA47").NumberFormat = "@"
$ws.Range("AA47").Value = '2023-08-27'
$ws.Range("AB47").Value = '00:00'
$ws.Range("AC47").Value = ""
$ws.Range("AD47").Value = $false
$ws.Range("AE47").Value = $false
$ws.Range("AF47").Value = ""
$ws.Range("AG47").Value = $false
$ws.Range("AH47").Value = ""
$ws.Range("AI47").Value = ""
$ws.Range("AJ47").Value = ""
$ws.Range("AK47").Value = ""
$ws.Range("AL47").Value = ""
$ws.Range("AM47").Value = ""
$ws.Range("AN47").Value = ""
$ws.Range("AO47").Value = ""
$ws.Range("AP47").Value = ""
$ws.Range("AQ47").Value = ""
$ws.Range("AR47").Value = ""
$ws.Range("AS47").Value = ""
$ws.Range("AT47").Value = ""
$ws.Range("AU47").Value = ""
$ws.Range("AV47").Value = ""
$ws.Range("AW47").Value = 'Jonas Nordenström'
$ws.Range("AX47").Value = 'Jonas Nordenström'
$ws.Range("AY47").Value = ""

# ---- Row 48 ----
$ws.Range("A48").Value = 111736370
$ws.Range("B48").Value = 56398
$ws.Range("C48").Value = 'Ovaliderad'
$ws.Range("D48").Value = 'NT'
$ws.Range("E48").Value = 100109
$ws.Range("F48").Value = 'Tretåig hackspett'
$ws.Range("G48").Value = 'Picoides tridactylus'
$ws.Range("H48").Value = '(Linnaeus, 1758)'
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""
$ws.Range("L48").Value = ""
$ws.Range("M48").Value = 'färska spår'
$ws.Range("N48").Value = ""
$ws.Range("O48").Value = ""
$ws.Range("P48").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q48").Value = 616327.1020967637
$ws.Range("R48").Value = 7268872.304318298
$ws.Range("S48").Value = 10
$ws.Range("T48").Value = 'Västerbotten'
$ws.Range("U48").Value = 'Sorsele'
$ws.Range("V48").Value = 'Lycksele lappmark'
$ws.Range("W48").Value = 'Sorsele'
$ws.Range("X48").Value = ""
$ws.Range("Y48").NumberFormat = "@"
$ws.Range("Y48").Value = '2023-08-27'
$ws.Range("Z48").Value = '00:00'
$ws.Range("AA48").NumberFormat = "@"
$ws.Range("AA48").Value = '2023-08-27'
$ws.Range("AB48").Value = '00:00'
$ws.Range("AC48").Value = ""
$ws.Range("AD48").Value = $false
$ws.Range("AE48").Value = $false
$ws.Range("AF48").Value = ""
$ws.Range("AG48").Value = $false
$ws.Range("AH48").Value = ""
$ws.Range("AI48").Value = ""
$ws.Range("AJ48").Value = ""
$ws.Range("AK48").Value = ""
$ws.Range("AL48").Value = ""
$ws.Range("AM48").Value = ""
$ws.Range("AN48").Value = ""
$ws.Range("AO48").Value = ""
$ws.Range("AP48").Value = ""
$ws.Range("AQ48").Value = ""
$ws.Range("AR48").Value = ""
$ws.Range("AS48").Value = ""
$ws.Range("AT48").Value = ""
$ws.Range("AU48").Value = ""
$ws.Range("AV48").Value = ""
$ws.Range("AW48").Value = 'Jonas Nordenström'
$ws.Range("AX48").Value = 'Jonas Nordenström'
$ws.Range("AY48").Value = ""

# ---- Row 49 ----
$ws.Range("A49").Value = 111736257
$ws.Range("B49").Value = 77515
$ws.Range("C49").Value = 'Ovaliderad'
$ws.Range("D49").Value = 'NT'
$ws.Range("E49").Value = 6425
$ws.Range("F49").Value = 'Garnlav'
$ws.Range("G49").Value = 'Alectoria sarmentosa'
$ws.Range("H49").Value = '(Ach.) Ach.'
$ws.Range("I49").Value = ""
$ws.Range("J49").Value = ""
$ws.Range("K49").Value = ""
$ws.Range("L49").Value = ""
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = ""
$ws.Range("O49").Value = ""
$ws.Range("P49").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q49").Value = 616308.8236423519
$ws.Range("R49").Value = 7268903.133137755
$ws.Range("S49").Value = 10
$ws.Range("T49").Value = 'Västerbotten'
$ws.Range("U49").Value = 'Sorsele'
$ws.Range("V49").Value = 'Lycksele lappmark'
$ws.Range("W49").Value = 'Sorsele'
$ws.Range("X49").Value = ""
$ws.Range("Y49").NumberFormat = "@"
$ws.Range("Y49").Value = '2023-08-27'
$ws.Range("Z49").Value = '00:00'
$ws.Range("AA49").NumberFormat = "@"
$ws.Range("AA49").Value = '2023-08-27'
$ws.Range("AB49").Value = '00:00'
$ws.Range("AC49").Value = ""
$ws.Range("AD49").Value = $false
$ws.Range("AE49").Value = $false
$ws.Range("AF49").Value = ""
$ws.Range("AG49").Value = $false
$ws.Range("AH49").Value = ""
$ws.Range("AI49").Value = ""
$ws.Range("AJ49").Value = ""
$ws.Range("AK49").Value = ""
$ws.Range("AL49").Value = ""
$ws.Range("AM49").Value = ""
$ws.Range("AN49").Value = ""
$ws.Range("AO49").Value = ""
$ws.Range("AP49").Value = ""
$ws.Range("AQ49").Value = ""
$ws.Range("AR49").Value = ""
$ws.Range("AS49").Value = ""
$ws.Range("AT49").Value = ""
$ws.Range("AU49").Value = ""
$ws.Range("AV49").Value = ""
$ws.Range("AW49").Value = 'Jonas Nordenström'
$ws.Range("AX49").Value = 'Jonas Nordenström'
$ws.Range("AY49").Value = ""

# ---- Row 51 ----
$ws.Range("A51").Value = 111736405
$ws.Range("B51").Value = 77515
$ws.Range("C51").Value = 'Ovaliderad'
$ws.Range("D51").Value = 'NT'
$ws.Range("E51").Value = 6425
$ws.Range("F51").Value = 'Garnlav'
$ws.Range("G51").Value = 'Alectoria sarmentosa'
$ws.Range("H51").Value = '(Ach.) Ach.'
$ws.Range("I51").Value = ""
$ws.Range("J51").Value = ""
$ws.Range("K51").Value = ""
$ws.Range("L51").Value = ""
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = ""
$ws.Range("O51").Value = ""
$ws.Range("P51").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q51").Value = 616333.1117616051
$ws.Range("R51").Value = 7268857.179896916
$ws.Range("S51").Value = 10
$ws.Range("T51").Value = 'Västerbotten'
$ws.Range("U51").Value = 'Sorsele'
$ws.Range("V51").Value = 'Lycksele lappmark'
$ws.Range("W51").Value = 'Sorsele'
$ws.Range("X51").Value = ""
$ws.Range("Y51").NumberFormat = "@"
$ws.Range("Y51").Value = '2023-08-27'
$ws.Range("Z51").Value = '00:00'
$ws.Range("AA51").NumberFormat = "@"
$ws.Range("AA51").Value = '2023-08-27'
$ws.Range("AB51").Value = '00:00'
$ws.Range("AC51").Value = ""
$ws.Range("AD51").Value = $false
$ws.Range("AE51").Value = $false
$ws.Range("AF51").Value = ""
$ws.Range("AG51").Value = $false
$ws.Range("AH51").Value = ""
$ws.Range("AI51").Value = ""
$ws.Range("AJ51").Value = ""
$ws.Range("AK51").Value = ""
$ws.Range("AL51").Value = ""
$ws.Range("AM51").Value = ""
$ws.Range("AN51").Value = ""
$ws.Range("AO51").Value = ""
$ws.Range("AP51").Value = ""
$ws.Range("AQ51").Value = ""
$ws.Range("AR51").Value = ""
$ws.Range("AS51").Value = ""
$ws.Range("AT51").Value = ""
$ws.Range("AU51").Value = ""
$ws.Range("AV51").Value = ""
$ws.Range("AW51").Value = 'Jonas Nordenström'
$ws.Range("AX51").Value = 'Jonas Nordenström'
$ws.Range("AY51").Value = ""

# ---- Row 52 ----
$ws.Range("A52").Value = 111736506
$ws.Range("B52").Value = 56398
$ws.Range("C52").Value = 'Ovaliderad'
$ws.Range("D52").Value = 'NT'
$ws.Range("E52").Value = 100109
$ws.Range("F52").Value = 'Tretåig hackspett'
$ws.Range("G52").Value = 'Picoides tridactylus'
$ws.Range("H52").Value = '(Linnaeus, 1758)'
$ws.Range("I52").Value = ""
$ws.Range("J52").Value = ""
$ws.Range("K52").Value = ""
$ws.Range("L52").Value = ""
$ws.Range("M52").Value = 'födosökande'
$ws.Range("N52").Value = ""
$ws.Range("O52").Value = ""
$ws.Range("P52").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q52").Value = 616358.6131022752
$ws.Range("R52").Value = 7268822.486957001
$ws.Range("S52").Value = 25
$ws.Range("T52").Value = 'Västerbotten'
$ws.Range("U52").Value = 'Sorsele'
$ws.Range("V52").Value = 'Lycksele lappmark'
$ws.Range("W52").Value = 'Sorsele'
$ws.Range("X52").Value = ""
$ws.Range("Y52").NumberFormat = "@"
$ws.Range("Y52").Value = '2023-08-27'
$ws.Range("Z52").Value = '00:00'
$ws.Range("AA52").NumberFormat = "@"
$ws.Range("AA52").Value = '2023-08-27'
$ws.Range("AB52").Value = '00:00'
$ws.Range("AC52").Value = 'Födosökande i gransumpskog. Observerades i över en timme'
$ws.Range("AD52").Value = $false
$ws.Range("AE52").Value = $false
$ws.Range("AF52").Value = ""
$ws.Range("AG52").Value = $false
$ws.Range("AH52").Value = ""
$ws.Range("AI52").Value = ""
$ws.Range("AJ52").Value = ""
$ws.Range("AK52").Value = ""
$ws.Range("AL52").Value = ""
$ws.Range("AM52").Value = ""
$ws.Range("AN52").Value = ""
$ws.Range("AO52").Value = ""
$ws.Range("AP52").Value = ""
$ws.Range("AQ52").Value = ""
$ws.Range("AR52").Value = ""
$ws.Range("AS52").Value = ""
$ws.Range("AT52").Value = ""
$ws.Range("AU52").Value = ""
$ws.Range("AV52").Value = ""
$ws.Range("AW52").Value = 'Jonas Nordenström'
$ws.Range("AX52").Value = 'Jonas Nordenström'
$ws.Range("AY52").Value = ""

# ---- Row 53 ----
$ws.Range("A53").Value = 111778248
$ws.Range("B53").Value = 56398
$ws.Range("C53").Value = 'Ovaliderad'
$ws.Range("D53").Value = 'NT'
$ws.Range("E53").Value = 100109
$ws.Range("F53").Value = 'Tretåig hackspett'
$ws.Range("G53").Value = 'Picoides tridactylus'
$ws.Range("H53").Value = '(Linnaeus, 1758)'
$ws.Range("I53").Value = ""
$ws.Range("J53").Value = ""
$ws.Range("K53").Value = ""
$ws.Range("L53").Value = ""
$ws.Range("M53").Value = 'färsk spillning'
$ws.Range("N53").Value = ""
$ws.Range("O53").Value = ""
$ws.Range("P53").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q53").Value = 616162.9874832245
$ws.Range("R53").Value = 7268630.281087617
$ws.Range("S53").Value = 10
$ws.Range("T53").Value = 'Västerbotten'
$ws.Range("U53").Value = 'Sorsele'
$ws.Range("V53").Value = 'Lycksele lappmark'
$ws.Range("W53").Value = 'Sorsele'
$ws.Range("X53").Value = ""
$ws.Range("Y53").NumberFormat = "@"
$ws.Range("Y53").Value = '2023-08-29'
$ws.Range("Z53").Value = '00:00'
$ws.Range("AA53").NumberFormat = "@"
$ws.Range("AA53").Value = '2023-08-29'
$ws.Range("AB53").Value = '00:00'
$ws.Range("AC53").Value = 'Skalad gran'
$ws.Range("AD53").Value = $false
$ws.Range("AE53").Value = $false
$ws.Range("AF53").Value = ""
$ws.Range("AG53").Value = $false
$ws.Range("AH53").Value = ""
$ws.Range("AI53").Value = ""
$ws.Range("AJ53").Value = ""
$ws.Range("AK53").Value = ""
$ws.Range("AL53").Value = ""
$ws.Range("AM53").Value = ""
$ws.Range("AN53").Value = ""
$ws.Range("AO53").Value = ""
$ws.Range("AP53").Value = ""
$ws.Range("AQ53").Value = ""
$ws.Range("AR53").Value = ""
$ws.Range("AS53").Value = ""
$ws.Range("AT53").Value = ""
$ws.Range("AU53").Value = ""
$ws.Range("AV53").Value = ""
$ws.Range("AW53").Value = 'Jonas Nordenström'
$ws.Range("AX53").Value = 'Jonas Nordenström'
$ws.Range("AY53").Value = ""

# ---- Row 54 ----
$ws.Range("A54").Value = 111777494
$ws.Range("B54").Value = 90854
$ws.Range("C54").Value = 'Ovaliderad'
$ws.Range("D54").Value = 'NT'
$ws.Range("E54").Value = 2079
$ws.Range("F54").Value = 'Nordtagging'
$ws.Range("G54").Value = 'Odonticium romellii'
$ws.Range("H54").Value = '(S.Lundell) Parmasto'
$ws.Range("I54").Value = ""
$ws.Range("J54").Value = ""
$ws.Range("K54").Value = ""
$ws.Range("L54").Value = ""
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = ""
$ws.Range("O54").Value = ""
$ws.Range("P54").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q54").Value = 616426.5202303537
$ws.Range("R54").Value = 7268746.301918368
$ws.Range("S54").Value = 10
$ws.Range("T54").Value = 'Västerbotten'
$ws.Range("U54").Value = 'Sorsele'
$ws.Range("V54").Value = 'Lycksele lappmark'
$ws.Range("W54").Value = 'Sorsele'
$ws.Range("X54").Value = ""
$ws.Range("Y54").NumberFormat = "@"
$ws.Range("Y54").Value = '2023-08-29'
$ws.Range("Z54").Value = '00:00'
$ws.Range("AA54").NumberFormat = "@"
$ws.Range("AA54").Value = '2023-08-29'
$ws.Range("AB54").Value = '00:00'
$ws.Range("AC54").Value = ""
$ws.Range("AD54").Value = $false
$ws.Range("AE54").Value = $false
$ws.Range("AF54").Value = ""
$ws.Range("AG54").Value = $false
$ws.Range("AH54").Value = ""
$ws.Range("AI54").Value = ""
$ws.Range("AJ54").Value = ""
$ws.Range("AK54").Value = ""
$ws.Range("AL54").Value = ""
$ws.Range("AM54").Value = ""
$ws.Range("AN54").Value = ""
$ws.Range("AO54").Value = ""
$ws.Range("AP54").Value = ""
$ws.Range("AQ54").Value = ""
$ws.Range("AR54").Value = ""
$ws.Range("AS54").Value = ""
$ws.Range("AT54").Value = ""
$ws.Range("AU54").Value = ""
$ws.Range("AV54").Value = ""
$ws.Range("AW54").Value = 'Jonas Nordenström'
$ws.Range("AX54").Value = 'Jonas Nordenström'
$ws.Range("AY54").Value = ""

# ---- Row 55 ----
$ws.Range("A55").Value = 111778163
$ws.Range("B55").Value = 56398
$ws.Range("C55").Value = 'Ovaliderad'
$ws.Range("D55").Value = 'NT'
$ws.Range("E55").Value = 100109
$ws.Range("F55").Value = 'Tretåig hackspett'
$ws.Range("G55").Value = 'Picoides tridactylus'
$ws.Range("H55").Value = '(Linnaeus, 1758)'
$ws.Range("I55").Value = ""
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = 'färska spår'
$ws.Range("N55").Value = ""
$ws.Range("O55").Value = ""
$ws.Range("P55").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q55").Value = 616207.2556492372
$ws.Range("R55").Value = 7268635.7870906
$ws.Range("S55").Value = 10
$ws.Range("T55").Value = 'Västerbotten'
$ws.Range("U55").Value = 'Sorsele'
$ws.Range("V55").Value = 'Lycksele lappmark'
$ws.Range("W55").Value = 'Sorsele'
$ws.Range("X55").Value = ""
$ws.Range("Y55").NumberFormat = "@"
$ws.Range("Y55").Value = '2023-08-29'
$ws.Range("Z55").Value = '00:00'
$ws.Range("AA55").NumberFormat = "@"
$ws.Range("AA55").Value = '2023-08-29'
$ws.Range("AB55").Value = '00:00'
$ws.Range("AC55").Value = 'Skalad gran'
$ws.Range("AD55").Value = $false
$ws.Range("AE55").Value = $false
$ws.Range("AF55").Value = ""
$ws.Range("AG55").Value = $false
$ws.Range("AH55").Value = ""
$ws.Range("AI55").Value = ""
$ws.Range("AJ55").Value = ""
$ws.Range("AK55").Value = ""
$ws.Range("AL55").Value = ""
$ws.Range("AM55").Value = ""
$ws.Range("AN55").Value = ""
$ws.Range("AO55").Value = ""
$ws.Range("AP55").Value = ""
$ws.Range("AQ55").Value = ""
$ws.Range("AR55").Value = ""
$ws.Range("AS55").Value = ""
$ws.Range("AT55").Value = ""
$ws.Range("AU55").Value = ""
$ws.Range("AV55").Value = ""
$ws.Range("AW55").Value = 'Jonas Nordenström'
$ws.Range("AX55").Value = 'Jonas Nordenström'
$ws.Range("AY55").Value = ""

# ---- Row 56 ----
$ws.Range("A56").Value = 111777331
$ws.Range("B56").Value = 89423
$ws.Range("C56").Value = 'Ovaliderad'
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 5432
$ws.Range("F56").Value = 'Granticka'
$ws.Range("G56").Value = 'Porodaedalea chrysoloma'
$ws.Range("H56").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("I56").Value = ""
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
$ws.Range("N56").Value = ""
$ws.Range("O56").Value = ""
$ws.Range("P56").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q56").Value = 616362.7639770868
$ws.Range("R56").Value = 7268822.653031595
$ws.Range("S56").Value = 10
$ws.Range("T56").Value = 'Västerbotten'
$ws.Range("U56").Value = 'Sorsele'
$ws.Range("V56").Value = 'Lycksele lappmark'
$ws.Range("W56").Value = 'Sorsele'
$ws.Range("X56").Value = ""
$ws.Range("Y56").NumberFormat = "@"
$ws.Range("Y56").Value = '2023-08-29'
$ws.Range("Z56").Value = '00:00'
$ws.Range("AA56").NumberFormat = "@"
$ws.Range("AA56").Value = '2023-08-29'
$ws.Range("AB56").Value = '00:00'
$ws.Range("AC56").Value = ""
$ws.Range("AD56").Value = $false
$ws.Range("AE56").Value = $false
$ws.Range("AF56").Value = ""
$ws.Range("AG56").Value = $false
$ws.Range("AH56").Value = ""
$ws.Range("AI56").Value = ""
$ws.Range("AJ56").Value = ""
$ws.Range("AK56").Value = ""
$ws.Range("AL56").Value = ""
$ws.Range("AM56").Value = ""
$ws.Range("AN56").Value = ""
$ws.Range("AO56").Value = ""
$ws.Range("AP56").Value = ""
$ws.Range("AQ56").Value = ""
$ws.Range("AR56").Value = ""
$ws.Range("AS56").Value = ""
$ws.Range("AT56").Value = ""
$ws.Range("AU56").Value = ""
$ws.Range("AV56").Value = ""
$ws.Range("AW56").Value = 'Jonas Nordenström'
$ws.Range("AX56").Value = 'Jonas Nordenström'
$ws.Range("AY56").Value = ""

# ---- Row 57 ----
$ws.Range("A57").Value = 111777380
$ws.Range("B57").Value = 56398
$ws.Range("C57").Value = 'Ovaliderad'
$ws.Range("D57").Value = 'NT'
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("I57").Value = ""
$ws.Range("J57").Value = ""
$ws.Range("K57").Value = ""
$ws.Range("L57").Value = ""
$ws.Range("M57").Value = 'färska spår'
$ws.Range("N57").Value = ""
$ws.Range("O57").Value = ""
$ws.Range("P57").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q57").Value = 616414.0528149965
$ws.Range("R57").Value = 7268860.418718725
$ws.Range("S57").Value = 10
$ws.Range("T57").Value = 'Västerbotten'
$ws.Range("U57").Value = 'Sorsele'
$ws.Range("V57").Value = 'Lycksele lappmark'
$ws.Range("W57").Value = 'Sorsele'
$ws.Range("X57").Value = ""
$ws.Range("Y57").NumberFormat = "@"
$ws.Range("Y57").Value = '2023-08-29'
$ws.Range("Z57").Value = '00:00'
$ws.Range("AA57").NumberFormat = "@"
$ws.Range("AA57").Value = '2023-08-29'
$ws.Range("AB57").Value = '00:00'
$ws.Range("AC57").Value = 'Skalade stammar'
$ws.Range("AD57").Value = $false
$ws.Range("AE57").Value = $false
$ws.Range("AF57").Value = ""
$ws.Range("AG57").Value = $false
$ws.Range("AH57").Value = ""
$ws.Range("AI57").Value = ""
$ws.Range("AJ57").Value = ""
$ws.Range("AK57").Value = ""
$ws.Range("AL57").Value = ""
$ws.Range("AM57").Value = ""
$ws.Range("AN57").Value = ""
$ws.Range("AO57").Value = ""
$ws.Range("AP57").Value = ""
$ws.Range("AQ57").Value = ""
$ws.Range("AR57").Value = ""
$ws.Range("AS57").Value = ""
$ws.Range("AT57").Value = ""
$ws.Range("AU57").Value = ""
$ws.Range("AV57").Value = ""
$ws.Range("AW57").Value = 'Jonas Nordenström'
$ws.Range("AX57").Value = 'Jonas Nordenström'
$ws.Range("AY57").Value = ""

# ---- Row 58 ----
$ws.Range("A58").Value = 111777447
$ws.Range("B58").Value = 89405
$ws.Range("C58").Value = 'Ovaliderad'
$ws.Range("D58").Value = 'NT'
$ws.Range("E58").Value = 1202
$ws.Range("F58").Value = 'Ullticka'
$ws.Range("G58").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H58").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I58").Value = ""
$ws.Range("J58").Value = ""
$ws.Range("K58").Value = ""
$ws.Range("L58").Value = ""
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = ""
$ws.Range("O58").Value = ""
$ws.Range("P58").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q58").Value = 616379.7321599644
$ws.Range("R58").Value = 7268803.814155157
$ws.Range("S58").Value = 10
$ws.Range("T58").Value = 'Västerbotten'
$ws.Range("U58").Value = 'Sorsele'
$ws.Range("V58").Value = 'Lycksele lappmark'
$ws.Range("W58").Value = 'Sorsele'
$ws.Range("X58").Value = ""
$ws.Range("Y58").NumberFormat = "@"
$ws.Range("Y58").Value = '2023-08-29'
$ws.Range("Z58").Value = '00:00'
$ws.Range("AA58").NumberFormat = "@"
$ws.Range("AA58").Value = '2023-08-29'
$ws.Range("AB58").Value = '00:00'
$ws.Range("AC58").Value = ""
$ws.Range("AD58").Value = $false
$ws.Range("AE58").Value = $false
$ws.Range("AF58").Value = ""
$ws.Range("AG58").Value = $false
$ws.Range("AH58").Value = ""
$ws.Range("AI58").Value = ""
$ws.Range("AJ58").Value = ""
$ws.Range("AK58").Value = ""
$ws.Range("AL58").Value = ""
$ws.Range("AM58").Value = ""
$ws.Range("AN58").Value = ""
$ws.Range("AO58").Value = ""
$ws.Range("AP58").Value = ""
$ws.Range("AQ58").Value = ""
$ws.Range("AR58").Value = ""
$ws.Range("AS58").Value = ""
$ws.Range("AT58").Value = ""
$ws.Range("AU58").Value = ""
$ws.Range("AV58").Value = ""
$ws.Range("AW58").Value = 'Jonas Nordenström'
$ws.Range("AX58").Value = 'Jonas Nordenström'
$ws.Range("AY58").Value = ""

# ---- Row 61 ----
$ws.Range("A61").Value = 111778005
$ws.Range("B61").Value = 89369
$ws.Range("C61").Value = 'Ovaliderad'
$ws.Range("D61").Value = 'LC'
$ws.Range("E61").Value = 5447
$ws.Range("F61").Value = 'Vedticka'
$ws.Range("G61").Value = 'Fuscoporia viticola'
$ws.Range("H61").Value = '(Schwein.) Murrill'
$ws.Range("I61").Value = ""
$ws.Range("J61").Value = ""
$ws.Range("K61").Value = ""
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = ""
$ws.Range("O61").Value = ""
$ws.Range("P61").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q61").Value = 616499.3130462242
$ws.Range("R61").Value = 7268610.508796399
$ws.Range("S61").Value = 10
$ws.Range("T61").Value = 'Västerbotten'
$ws.Range("U61").Value = 'Sorsele'
$ws.Range("V61").Value = 'Lycksele lappmark'
$ws.Range("W61").Value = 'Sorsele'
$ws.Range("X61").Value = ""
$ws.Range("Y61").NumberFormat = "@"
$ws.Range("Y61").Value = '2023-08-29'
$ws.Range("Z61").Value = '00:00'
$ws.Range("AA61").NumberFormat = "@"
$ws.Range("AA61").Value = '2023-08-29'
$ws.Range("AB61").Value = '00:00'
$ws.Range("AC61").Value = ""
$ws.Range("AD61").Value = $false
$ws.Range("AE61").Value = $false
$ws.Range("AF61").Value = ""
$ws.Range("AG61").Value = $false
$ws.Range("AH61").Value = ""
$ws.Range("AI61").Value = ""
$ws.Range("AJ61").Value = ""
$ws.Range("AK61").Value = ""
$ws.Range("AL61").Value = ""
$ws.Range("AM61").Value = ""
$ws.Range("AN61").Value = ""
$ws.Range("AO61").Value = ""
$ws.Range("AP61").Value = ""
$ws.Range("AQ61").Value = ""
$ws.Range("AR61").Value = ""
$ws.Range("AS61").Value = ""
$ws.Range("AT61").Value = ""
$ws.Range("AU61").Value = ""
$ws.Range("AV61").Value = ""
$ws.Range("AW61").Value = 'Jonas Nordenström'
$ws.Range("AX61").Value = 'Jonas Nordenström'
$ws.Range("AY61").Value = ""

# ---- Row 62 ----
$ws.Range("A62").Value = 111778126
$ws.Range("B62").Value = 89405
$ws.Range("C62").Value = 'Ovaliderad'
$ws.Range("D62").Value = 'NT'
$ws.Range("E62").Value = 1202
$ws.Range("F62").Value = 'Ullticka'
$ws.Range("G62").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H62").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I62").Value = ""
$ws.Range("J62").Value = ""
$ws.Range("K62").Value = ""
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("O62").Value = ""
$ws.Range("P62").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q62").Value = 616202.3044715263
$ws.Range("R62").Value = 7268603.611313918
$ws.Range("S62").Value = 10
$ws.Range("T62").Value = 'Västerbotten'
$ws.Range("U62").Value = 'Sorsele'
$ws.Range("V62").Value = 'Lycksele lappmark'
$ws.Range("W62").Value = 'Sorsele'
$ws.Range("X62").Value = ""
$ws.Range("Y62").NumberFormat = "@"
$ws.Range("Y62").Value = '2023-08-29'
$ws.Range("Z62").Value = '00:00'
$ws.Range("AA62").NumberFormat = "@"
$ws.Range("AA62").Value = '2023-08-29'
$ws.Range("AB62").Value = '00:00'
$ws.Range("AC62").Value = ""
$ws.Range("AD62").Value = $false
$ws.Range("AE62").Value = $false
$ws.Range("AF62").Value = ""
$ws.Range("AG62").Value = $false
$ws.Range("AH62").Value = ""
$ws.Range("AI62").Value = ""
$ws.Range("AJ62").Value = ""
$ws.Range("AK62").Value = ""
$ws.Range("AL62").Value = ""
$ws.Range("AM62").Value = ""
$ws.Range("AN62").Value = ""
$ws.Range("AO62").Value = ""
$ws.Range("AP62").Value = ""
$ws.Range("AQ62").Value = ""
$ws.Range("AR62").Value = ""
$ws.Range("AS62").Value = ""
$ws.Range("AT62").Value = ""
$ws.Range("AU62").Value = ""
$ws.Range("AV62").Value = ""
$ws.Range("AW62").Value = 'Jonas Nordenström'
$ws.Range("AX62").Value = 'Jonas Nordenström'
$ws.Range("AY62").Value = ""

# ---- Row 63 ----
$ws.Range("A63").Value = 111777411
$ws.Range("B63").Value = 56398
$ws.Range("C63").Value = 'Ovaliderad'
$ws.Range("D63").Value = 'NT'
$ws.Range("E63").Value = 100109
$ws.Range("F63").Value = 'Tretåig hackspett'
$ws.Range("G63").Value = 'Picoides tridactylus'
$ws.Range("H63").Value = '(Linnaeus, 1758)'
$ws.Range("I63").Value = ""
$ws.Range("J63").Value = ""
$ws.Range("K63").Value = ""
$ws.Range("L63").Value = ""
$ws.Range("M63").Value = 'färska spår'
$ws.Range("N63").Value = ""
$ws.Range("O63").Value = ""
$ws.Range("P63").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q63").Value = 616367.7277224116
$ws.Range("R63").Value = 7268802.503264537
$ws.Range("S63").Value = 10
$ws.Range("T63").Value = 'Västerbotten'
$ws.Range("U63").Value = 'Sorsele'
$ws.Range("V63").Value = 'Lycksele lappmark'
$ws.Range("W63").Value = 'Sorsele'
$ws.Range("X63").Value = ""
$ws.Range("Y63").NumberFormat = "@"
$ws.Range("Y63").Value = '2023-08-29'
$ws.Range("Z63").Value = '00:00'
$ws.Range("AA63").NumberFormat = "@"
$ws.Range("AA63").Value = '2023-08-29'
$ws.Range("AB63").Value = '00:00'
$ws.Range("AC63").Value = 'Skalade granstammar'
$ws.Range("AD63").Value = $false
$ws.Range("AE63").Value = $false
$ws.Range("AF63").Value = ""
$ws.Range("AG63").Value = $false
$ws.Range("AH63").Value = ""
$ws.Range("AI63").Value = ""
$ws.Range("AJ63").Value = ""
$ws.Range("AK63").Value = ""
$ws.Range("AL63").Value = ""
$ws.Range("AM63").Value = ""
$ws.Range("AN63").Value = ""
$ws.Range("AO63").Value = ""
$ws.Range("AP63").Value = ""
$ws.Range("AQ63").Value = ""
$ws.Range("AR63").Value = ""
$ws.Range("AS63").Value = ""
$ws.Range("AT63").Value = ""
$ws.Range("AU63").Value = ""
$ws.Range("AV63").Value = ""
$ws.Range("AW63").Value = 'Jonas Nordenström'
$ws.Range("AX63").Value = 'Jonas Nordenström'
$ws.Range("AY63").Value = ""

# ---- Row 64 ----
$ws.Range("A64").Value = 111777499
$ws.Range("B64").Value = 78107
$ws.Range("C64").Value = 'Ovaliderad'
$ws.Range("D64").Value = 'NT'
$ws.Range("E64").Value = 6453
$ws.Range("F64").Value = 'Vedskivlav'
$ws.Range("G64").Value = 'Hertelidea botryosa'
$ws.Range("H64").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("I64").Value = ""
$ws.Range("J64").Value = ""
$ws.Range("K64").Value = ""
$ws.Range("L64").Value = ""
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = ""
$ws.Range("O64").Value = ""
$ws.Range("P64").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q64").Value = 616426.5202303537
$ws.Range("R64").Value = 7268746.301918368
$ws.Range("S64").Value = 10
$ws.Range("T64").Value = 'Västerbotten'
$ws.Range("U64").Value = 'Sorsele'
$ws.Range("V64").Value = 'Lycksele lappmark'
$ws.Range("W64").Value = 'Sorsele'
$ws.Range("X64").Value = ""
$ws.Range("Y64").NumberFormat = "@"
$ws.Range("Y64").Value = '2023-08-29'
$ws.Range("Z64").Value = '00:00'
$ws.Range("AA64").NumberFormat = "@"
$ws.Range("AA64").Value = '2023-08-29'
$ws.Range("AB64").Value = '00:00'
$ws.Range("AC64").Value = ""
$ws.Range("AD64").Value = $false
$ws.Range("AE64").Value = $false
$ws.Range("AF64").Value = ""
$ws.Range("AG64").Value = $false
$ws.Range("AH64").Value = ""
$ws.Range("AI64").Value = ""
$ws.Range("AJ64").Value = ""
$ws.Range("AK64").Value = ""
$ws.Range("AL64").Value = ""
$ws.Range("AM64").Value = ""
$ws.Range("AN64").Value = ""
$ws.Range("AO64").Value = ""
$ws.Range("AP64").Value = ""
$ws.Range("AQ64").Value = ""
$ws.Range("AR64").Value = ""
$ws.Range("AS64").Value = ""
$ws.Range("AT64").Value = ""
$ws.Range("AU64").Value = ""
$ws.Range("AV64").Value = ""
$ws.Range("AW64").Value = 'Jonas Nordenström'
$ws.Range("AX64").Value = 'Jonas Nordenström'
$ws.Range("AY64").Value = ""

# ---- Row 65 ----
$ws.Range("A65").Value = 111777491
$ws.Range("B65").Value = 56398
$ws.Range("C65").Value = 'Ovaliderad'
$ws.Range("D65").Value = 'NT'
$ws.Range("E65").Value = 100109
$ws.Range("F65").Value = 'Tretåig hackspett'
$ws.Range("G65").Value = 'Picoides tridactylus'
$ws.Range("H65").Value = '(Linnaeus, 1758)'
$ws.Range("I65").Value = ""
$ws.Range("J65").Value = ""
$ws.Range("K65").Value = ""
$ws.Range("L65").Value = ""
$ws.Range("M65").Value = 'färska spår'
$ws.Range("N65").Value = ""
$ws.Range("O65").Value = ""
$ws.Range("P65").Value = 'Rankbäcken, Ly lm'
$ws.Range("Q65").Value = 616426.5202303537
$ws.Range("R65").Value = 7268746.301918368
$ws.Range("S65").Value = 10
$ws.Range("T65").Value = 'Västerbotten'
$ws.Range("U65").Value = 'Sorsele'
$ws.Range("V65").Value = 'Lycksele lappmark'
$ws.Range("W65").Value = 'Sorsele'
$ws.Range("X65").Value = ""
$ws.Range("Y65").NumberFormat = "@"
$ws.Range("Y65").Value = '2023-08-29'
$ws.Range("Z65").Value = '00:00'
$ws.Range("AA65").NumberFormat = "@"
$ws.Range("AA65").Value = '2023-08-29'
$ws.Range("AB65").Value = '00:00'
$ws.Range("AC65").Value = 'Skalade granstammar'
$ws.Range("AD65").Value = $false
$ws.Range("AE65").Value = $false
$ws.Range("AF65").Value = ""
$ws.Range("AG65").Value = $false
$ws.Range("AH65").Value = ""
$ws.Range("AI65").Value = ""
$ws.Range("AJ65").Value = ""
$ws.Range("AK65").Value = ""
$ws.Range("AL65").Value = ""
$ws.Range("AM65").Value = ""
$ws.Range("AN65").Value = ""
$ws.Range("AO65").Value = ""
$ws.Range("AP65").Value = ""
$ws.Range("AQ65").Value = ""
$ws.Range("AR65").Value = ""
$ws.Range("AS65").Value = ""
$ws.Range("AT65").Value = ""
$ws.Range("AU65").Value = ""
$ws.Range("AV65").Value = ""
$ws.Range("AW65").Value = 'Jonas Nordenström'
$ws.Range("AX65").Value = 'Jonas Nordenström'
$ws.Range("AY65").Value = ""
